$d = $word.ActiveDocument

function Replace-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $oldText = $r.Text
    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Paragraph 2: "En primer lugar..." -> fix ScriptableObjects wording / Sting -> String
Replace-ParagraphText 2 "En primer lugar, pensé como estructurar los números que van saliendo por la pantalla. Al utilizar una estética concreta he usado ScriptableObjects para las cartas, que tienen tres componentes a guardar: el String del número, el int del número y el Sprite de la carta."

# Paragraph 4: "Esta parte..." -> new leading sentence, Sting -> String
Replace-ParagraphText 4 "En caso de que no se utilizaran ScriptableObjects esta parte se podría hacer de otra manera, ya que se podría diseñar una función que se encargara de traducir el número X a String. Por comodidad se ha estructurado con los ScriptableObjects donde es muy fácil gestionarlo todo con sus Sprites correspondientes."

# Paragraph 5: "La estética escogida..." -> merge runs only, same text
Replace-ParagraphText 5 "La estética escogida se ha decido para dar un enfoque al ejercicio y hacerlo más visual y divertido."

# Paragraph 6: "También he decidido..." -> "como" -> "cómo"
Replace-ParagraphText 6 "También he decidido que el juego se gestionaría entero des del GameController, allí se hacen todas las llamadas en orden, y después el numero aleatorio que va saliendo y las opciones tienen otros scripts para gestionar su animación de salida y cómo reaccionan ante las interacciones."

# Paragraph 7: "Hay un script..." becomes the first of five paragraphs about the scripts.
# The original paragraph node becomes the new "El script de Opciones..." paragraph,
# and four new paragraphs are inserted after it.
Replace-ParagraphText 7 "El script de Opciones se encarga de gestionar el clic en los botones, des de allí se comprueba si ha acertado o no y se modifica el color del botón. También se gestionan las animaciones de salida. Los botones se activan solamente cuando ya están en posición para que así no se cliquen por error cuando están saliendo."

$d.Paragraphs(7).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(8).Range.Text = "El script de los números aleatorios es más sencillo, y solo controla el fadeIn o fadeOut de los números. "

$d.Paragraphs(8).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(9).Range.Text = "Hay un script que se encarga de modificar los parámetros en pantalla de los aciertos y errores, en este caso solo son dos, pero en el caso que la interfaz fuese más completa se encargaría de gestionar todos los otros parámetros."

$d.Paragraphs(9).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(10).Range.Text = "Los scripts se han fraccionado para poder gestionar cosas distintas des de cada uno, así la depuración del código es más sencilla y está todo más ordenado."

$d.Paragraphs(10).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(11).Range.Text = "Igualmente, siempre hay el script principal des de donde se va controlando el flujo. Gracias a esto podemos saber que va pasando en cada momento y des de allí buscar donde están los problemas."

Write-Output "done"
